$wb = $excel.ActiveWorkbook

# --- Sheet1 ("CCK81_noCTRL_meas" measurements sheet) ---------------------
# Rows 45:87 were a stray leftover - just an index column (44..86) running
# past the real data block (A1:N44). Remove them so the sheet's data range
# matches the other two sheets (A1:N44).
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows("45:87").Delete()

# Restore Sheet1 as the active/selected tab (it had drifted to Sheet3) and
# reproduce the saved cursor position/selection (C68, scrolled to row 24).
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 24
$ws1.Range("C68").Select() | Out-Null
